$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.134.59"
$ws.Range("E2").Value = "  +0.00%  "
$ws.Range("D3").Value = "'1.654.46"
$ws.Range("E3").Value = "  -0.21%  "
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("D5").Value = "'218.53"
$ws.Range("E5").Value = "  -0.14%  "
$ws.Range("D6").Value = "'0.5239"
$ws.Range("E6").Value = "  -0.22%  "
$ws.Range("E7").Value = "  -0.19%  "
$ws.Range("D8").Value = "'0.2655"
$ws.Range("E8").Value = "  +1.29%  "
$ws.Range("D9").Value = "'0.06354"
$ws.Range("E10").Value = "  -0.46%  "
$ws.Range("D11").Value = "'0.07700"
$ws.Range("E11").Value = "  -1.33%  "
$ws.Range("D12").Value = "'4.636"
$ws.Range("E12").Value = "  +3.18%  "
$ws.Range("D13").Value = "'1.699.97"
$ws.Range("E13").Value = "  +2.49%  "
$ws.Range("D14").Value = "'1.882.44"
$ws.Range("E14").Value = "  -0.16%  "
$ws.Range("D15").Value = "'0.5615"
$ws.Range("E15").Value = "  +0.96%  "
$ws.Range("D16").Value = "'0.0₅8180"
$ws.Range("E16").Value = "  +1.80%  "
$ws.Range("E17").Value = "  +0.60%  "
$ws.Range("D18").Value = "'26.123.83"
$ws.Range("E19").Value = "  -0.22%  "
$ws.Range("D20").Value = "'4.651"
$ws.Range("E20").Value = "  +0.14%  "
$ws.Range("D21").Value = "'10.50"
$ws.Range("E21").Value = "  +3.89%  "
$ws.Range("D22").Value = "'192.35"
$ws.Range("E22").Value = "  -1.58%  "
$ws.Range("D23").Value = "'5.955"
$ws.Range("E23").Value = "  -0.16%  "
$ws.Range("E24").Value = "  -0.18%  "
$ws.Range("D25").Value = "'144.82"
$ws.Range("E25").Value = "  -1.70%  "
$ws.Range("E26").Value = "  -0.98%  "
$ws.Range("D27").Value = "'7.265"
$ws.Range("E27").Value = "  +1.17%  "
$ws.Range("E28").Value = "  +0.01%  "
$ws.Range("D29").Value = "'1.513"
$ws.Range("E29").Value = "  +1.15%  "
$ws.Range("D30").Value = "'0.05450"
$ws.Range("E30").Value = "  -4.65%  "
$ws.Range("E31").Value = "  -0.08%  "
$ws.Range("E32").Value = "  -0.78%  "
$ws.Range("E33").Value = "  +0.44%  "
$ws.Range("E34").Value = "  -1.73%  "
$ws.Range("D35").Value = "'0.9513"
$ws.Range("E35").Value = "  -0.16%  "
$ws.Range("D36").Value = "'2.780"
$ws.Range("E36").Value = "  -0.88%  "
$ws.Range("D37").Value = "'2.402"
$ws.Range("E37").Value = "  -0.61%  "
$ws.Range("D38").Value = "'0.5679"
$ws.Range("E38").Value = "  -0.54%  "
$ws.Range("D39").Value = "'0.01582"
$ws.Range("E39").Value = "  -0.96%  "
$ws.Range("D40").Value = "'5.865"
$ws.Range("E40").Value = "  -1.34%  "
$ws.Range("E41").Value = "  -0.17%  "
$ws.Range("D42").Value = "'0.8336"
$ws.Range("E42").Value = "  -1.54%  "
$ws.Range("D43").Value = "'1.027.68"
$ws.Range("E43").Value = "  -3.57%  "
$ws.Range("D44").Value = "'101.19"
$ws.Range("E44").Value = "  -2.28%  "
$ws.Range("D45").Value = "'1.793.22"
$ws.Range("E45").Value = "  -0.13%  "
$ws.Range("D46").Value = "'57.72"
$ws.Range("E46").Value = "  -0.21%  "
$ws.Range("D47").Value = "'0.9994"
$ws.Range("E47").Value = "  -0.80%  "
$ws.Range("E48").Value = "  -1.36%  "
$ws.Range("D49").Value = "'7.994"
$ws.Range("E49").Value = "  -0.31%  "
$ws.Range("E50").Value = "  -2.11%  "
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").Value = "'0.09714"
$ws.Range("E51").Value = "  +2.83%  "
